# BehaviorScenario_HouseholdComposition.xlsx
# "added a new type of household and added scenarios"
#
# The sheet holds household-composition scenarios:
#   col A = ID_HouseholdType, col B = ID_PersonType, col C = unit ("count"), col D = value
#
# This script:
#   1) tweaks the scenario counts for the existing household types 3 and 4
#   2) appends four new rows describing a brand-new household type (5),
#      one row per person type 1-4
#   3) leaves the selection where the author left it when they saved

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- revised scenario counts for household types 3 & 4 ---
$ws.Range("D10").Value = 2
$ws.Range("D11").Value = 0
$ws.Range("D12").Value = 1

$ws.Range("D14").Value = 1
$ws.Range("D15").Value = 1
$ws.Range("D16").Value = 2
$ws.Range("D17").Value = 0

# --- new household type 5 scenarios ---
$ws.Range("A18").Value = 5
$ws.Range("B18").Value = 1
$ws.Range("C18").Value = "count"
$ws.Range("D18").Value = 0

$ws.Range("A19").Value = 5
$ws.Range("B19").Value = 2
$ws.Range("C19").Value = "count"
$ws.Range("D19").Value = 0

$ws.Range("A20").Value = 5
$ws.Range("B20").Value = 3
$ws.Range("C20").Value = "count"
$ws.Range("D20").Value = 0

$ws.Range("A21").Value = 5
$ws.Range("B21").Value = 4
$ws.Range("C21").Value = "count"
$ws.Range("D21").Value = 2

# --- leave the cursor parked where the saved file shows it ---
$ws.Range("A8").Select()
